# Update the cryptocurrency price/volume table to reflect the latest
# scrape (commit: "Updated cryptos list ... with GitHub Actions").
#
# Column D values look numeric (e.g. "29.188.93", "0.07761") but are
# stored as plain text in the sheet (OOXML inlineStr, not <v> numbers) —
# some even use '.' as a thousands separator in a way Excel would
# otherwise parse as a number or reformat (dropping trailing zeros).
# Force those cells to Text before assigning so Excel keeps the exact
# string, then restore the "Normal" style so no stray number-format
# style is left behind on the cell.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "29.188.93"
$ws.Range("E2").Value = "  -0.65%  "
Set-TextValue $ws.Range("D3") "1.857.30"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D5") "241.86"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D6") "0.7015"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.3106"
$ws.Range("E8").Value = "  -0.58%  "
Set-TextValue $ws.Range("D9") "0.07761"
$ws.Range("E9").Value = "  -3.56%  "
Set-TextValue $ws.Range("D10") "24.11"
$ws.Range("E10").Value = "  -4.63%  "
Set-TextValue $ws.Range("D11") "0.07985"
$ws.Range("E11").Value = "  -4.14%  "
Set-TextValue $ws.Range("D12") "1.854.09"
$ws.Range("E12").Value = "  -1.96%  "
Set-TextValue $ws.Range("D13") "5.167"
$ws.Range("E13").Value = "  -1.48%  "
Set-TextValue $ws.Range("D14") "93.28"
$ws.Range("E14").Value = "  -0.38%  "
Set-TextValue $ws.Range("D15") "0.6944"
$ws.Range("E15").Value = "  -3.48%  "
Set-TextValue $ws.Range("D16") "6.339"
$ws.Range("E16").Value = "  +0.21%  "
Set-TextValue $ws.Range("D17") "29.168.16"
$ws.Range("E17").Value = "  -0.75%  "
Set-TextValue $ws.Range("D18") "0.000008267"
$ws.Range("E18").Value = "  -3.22%  "
Set-TextValue $ws.Range("D19") "250.09"
$ws.Range("E19").Value = "  +3.39%  "
Set-TextValue $ws.Range("D20") "2.112.64"
$ws.Range("E20").Value = "  -1.12%  "
Set-TextValue $ws.Range("D21") "13.05"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  -0.01%  "
Set-TextValue $ws.Range("D23") "7.493"
$ws.Range("E23").Value = "  -4.52%  "
Set-TextValue $ws.Range("D24") "1.002"
$ws.Range("E24").Value = "  +0.09%  "
Set-TextValue $ws.Range("D25") "0.1545"
$ws.Range("E25").Value = "  -2.60%  "
Set-TextValue $ws.Range("D26") "8.949"
$ws.Range("E26").Value = "  -1.24%  "
Set-TextValue $ws.Range("D27") "159.17"
$ws.Range("E27").Value = "  -2.94%  "
Set-TextValue $ws.Range("D28") "18.73"
$ws.Range("E28").Value = "  +0.75%  "
Set-TextValue $ws.Range("D29") "1.493"
$ws.Range("E29").Value = "  -0.96%  "
Set-TextValue $ws.Range("D30") "4.279"
$ws.Range("E30").Value = "  -3.15%  "
Set-TextValue $ws.Range("D31") "4.247"
$ws.Range("E31").Value = "  -2.19%  "
Set-TextValue $ws.Range("D32") "1.209"
$ws.Range("E32").Value = "  +0.85%  "
Set-TextValue $ws.Range("D33") "0.05241"
$ws.Range("E33").Value = "  -2.38%  "
Set-TextValue $ws.Range("D34") "1.872"
$ws.Range("E34").Value = "  -3.93%  "
Set-TextValue $ws.Range("D35") "0.7410"
$ws.Range("E35").Value = "  -1.06%  "
Set-TextValue $ws.Range("D36") "1.152"
$ws.Range("E36").Value = "  -2.64%  "
Set-TextValue $ws.Range("D37") "2.713"
$ws.Range("E37").Value = "  +0.62%  "
Set-TextValue $ws.Range("D38") "0.01860"
$ws.Range("E38").Value = "  -1.42%  "
Set-TextValue $ws.Range("D39") "1.242.75"
$ws.Range("E39").Value = "  -3.54%  "
Set-TextValue $ws.Range("D40") "2.733"
$ws.Range("E40").Value = "  -0.51%  "
Set-TextValue $ws.Range("D41") "6.206"
$ws.Range("E41").Value = "  -6.05%  "
Set-TextValue $ws.Range("D42") "110.63"
$ws.Range("E42").Value = "  -1.34%  "
Set-TextValue $ws.Range("D43") "0.8934"
$ws.Range("E43").Value = "  -2.85%  "
Set-TextValue $ws.Range("D44") "70.71"
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +0.55%  "
Set-TextValue $ws.Range("D47") "2.010.66"
$ws.Range("E47").Value = "  -1.35%  "
Set-TextValue $ws.Range("D48") "0.5184"
$ws.Range("E48").Value = "  -0.71%  "
Set-TextValue $ws.Range("D49") "1.774"
$ws.Range("E49").Value = "  -1.76%  "
Set-TextValue $ws.Range("D50") "9.387"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D51") "1.004"
$ws.Range("E51").Value = "  +0.05%  "
